$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "?" placeholders in column B for IOT / BlockChain rows
$ws.Range("B7").Value = "20k"
$ws.Range("B8").Value = "35k"

# Insert a row above the old "WorkShop" row (currently row 10) so new
# course rows (PowerBI, Tableau, IBM Watson Cloud) can be added, then the
# WorkShop row itself is moved down one more row (to row 12) and a blank
# row (13) is left before the Surya/Dinesh experience rows.
$ws.Rows("9:9").Insert()

# New course/amount rows
$ws.Range("A9").Value = "PowerBI"
$ws.Range("B9").Value = "15k"

$ws.Range("A10").Value = "Tableau"
$ws.Range("B10").Value = "15k"

$ws.Range("A11").Value = "IBM Watson Cloud"
$ws.Range("B11").Value = "30k"

# WorkShop row (was row 10, now pushed down to row 12)
$ws.Range("A12").Value = "WorkShop"
$ws.Range("B12").Value = "15k/Day"

# Move the Surya / Dinesh experience rows down one, from 12:13 to 14:15,
# leaving row 13 blank.
$ws.Range("A13").Value = $null
$ws.Range("C13").Value = $null

$ws.Range("A14").Value = "Surya"
$ws.Range("C14").Value = "4yrs"

$ws.Range("A15").Value = "Dinesh"
$ws.Range("C15").Value = "4yrs"

# Update the selection to match the author's saved cursor position
$ws.Range("F3").Select()
